$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 11.0429782294745
$ws.Range("C2").Value = 8.908426769806431
$ws.Range("D2").Value = 6.038570938916996
$ws.Range("E2").Value = 12.59899722849428
$ws.Range("F2").Value = 28.63153133745559
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 25.91616189452973
$ws.Range("K2").Value = 8.510001872591852
$ws.Range("L2").Value = 10.18228369431872
$ws.Range("M2").Value = 13.60919536096916
$ws.Range("O2").Value = 25.73360879268883
# Row 3
$ws.Range("B3").Value = 10.79659435792805
$ws.Range("C3").Value = 8.874182159500899
$ws.Range("D3").Value = 5.990179931017614
$ws.Range("E3").Value = 12.62948369414344
$ws.Range("F3").Value = 28.6883234186766
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 26.01196085566676
$ws.Range("K3").Value = 8.322947988888117
$ws.Range("L3").Value = 10.19043898025124
$ws.Range("M3").Value = 13.57305300077464
$ws.Range("O3").Value = 25.81635565935115
# Row 4
$ws.Range("B4").Value = 10.64411892195413
$ws.Range("C4").Value = 8.852898917257686
$ws.Range("D4").Value = 5.959787623884229
$ws.Range("E4").Value = 12.6498422482121
$ws.Range("F4").Value = 28.72990042668
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 26.07520621252328
$ws.Range("K4").Value = 8.207052698407255
$ws.Range("L4").Value = 10.19683424986208
$ws.Range("M4").Value = 13.55267232263873
$ws.Range("O4").Value = 25.87207355515976
# Row 5
$ws.Range("B5").Value = 10.58177202474375
$ws.Range("C5").Value = 8.844163446659428
$ws.Range("D5").Value = 5.947236017613102
$ws.Range("E5").Value = 12.65855117777661
$ws.Range("F5").Value = 28.74852682541429
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 26.10209152440307
$ws.Range("K5").Value = 8.159627079494939
$ws.Range("L5").Value = 10.19978990835933
$ws.Range("M5").Value = 13.54482817205532
$ws.Range("O5").Value = 25.89601247502313
# Row 6
$ws.Range("B6").Value = 10.57140921805473
$ws.Range("C6").Value = 8.842709225808644
$ws.Range("D6").Value = 5.945141909668408
$ws.Range("E6").Value = 12.6600222216185
$ws.Range("F6").Value = 28.75172130437073
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 26.10662298444701
$ws.Range("K6").Value = 8.151742118363289
$ws.Range("L6").Value = 10.20030181860291
$ws.Range("M6").Value = 13.54355367049512
$ws.Range("O6").Value = 25.9000619658892
# Row 7
$ws.Range("B7").Value = 10.64327882869904
$ws.Range("C7").Value = 8.8527813572758
$ws.Range("D7").Value = 5.959619016684587
$ws.Range("E7").Value = 12.64995802862966
$ws.Range("F7").Value = 28.73014481691914
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 26.07556429382846
$ws.Range("K7").Value = 8.20641381133022
$ws.Range("L7").Value = 10.19687269499216
$ws.Range("M7").Value = 13.55256465916662
$ws.Range("O7").Value = 25.87239141172261
# Row 8
$ws.Range("B8").Value = 10.9583251289025
$ws.Range("C8").Value = 8.896673885170356
$ws.Range("D8").Value = 6.022029543047568
$ws.Range("E8").Value = 12.60916883306034
$ws.Range("F8").Value = 28.6497196790944
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 25.94827494846869
$ws.Range("K8").Value = 8.44576052850198
$ws.Range("L8").Value = 10.1848079210054
$ws.Range("M8").Value = 13.59636119950369
$ws.Range("O8").Value = 25.76111976388545
# Row 9
$ws.Range("B9").Value = 11.56293756808465
$ws.Range("C9").Value = 8.980612265258369
$ws.Range("D9").Value = 6.138827366896424
$ws.Range("E9").Value = 12.54217726433015
$ws.Range("F9").Value = 28.54533338721726
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 25.73377701991349
$ws.Range("K9").Value = 8.904112479171689
$ws.Range("L9").Value = 10.1721357900633
$ws.Range("M9").Value = 13.69634651108438
$ws.Range("O9").Value = 25.58194735001238
# Row 10
$ws.Range("B10").Value = 11.99441598483723
$ws.Range("C10").Value = 9.040853725767503
$ws.Range("D10").Value = 6.220979651241077
$ws.Range("E10").Value = 12.5008609835164
$ws.Range("F10").Value = 28.50128037746827
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 25.59760700412052
$ws.Range("K10").Value = 9.230710984420549
$ws.Range("L10").Value = 10.16948654731713
$ws.Range("M10").Value = 13.77802261617256
$ws.Range("O10").Value = 25.47418913155659
# Row 11
$ws.Range("B11").Value = 12.18704041500562
$ws.Range("C11").Value = 9.067922028094323
$ws.Range("D11").Value = 6.257506367075222
$ws.Range("E11").Value = 12.48377663380294
$ws.Range("F11").Value = 28.48834414320976
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 25.54031461227218
$ws.Range("K11").Value = 9.376421774492224
$ws.Range("L11").Value = 10.16971840152535
$ws.Range("M11").Value = 13.81687589485274
$ws.Range("O11").Value = 25.43036949716316
# Row 12
$ws.Range("B12").Value = 12.25938774354644
$ws.Range("C12").Value = 9.078121432154086
$ws.Range("D12").Value = 6.271212279741426
$ws.Range("E12").Value = 12.47755286245497
$ws.Range("F12").Value = 28.4844674642215
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 25.51928884190477
$ws.Range("K12").Value = 9.431136848611565
$ws.Range("L12").Value = 10.17001190338495
$ws.Range("M12").Value = 13.8318249267502
$ws.Range("O12").Value = 25.41452524178099
# Row 13
$ws.Range("B13").Value = 12.24383398330527
$ws.Range("C13").Value = 9.075927105726887
$ws.Range("D13").Value = 6.268266138781121
$ws.Range("E13").Value = 12.47888234002377
$ws.Range("F13").Value = 28.48525692147336
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 25.52378732582824
$ws.Range("K13").Value = 9.419374320355347
$ws.Range("L13").Value = 10.16993955954629
$ws.Range("M13").Value = 13.82859501478148
$ws.Range("O13").Value = 25.41790423754911
# Row 14
$ws.Range("B14").Value = 12.19300478520255
$ws.Range("C14").Value = 9.068762171017905
$ws.Range("D14").Value = 6.258636510566923
$ws.Range("E14").Value = 12.48325967784959
$ws.Range("F14").Value = 28.48800472353018
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 25.53857138559192
$ws.Range("K14").Value = 9.380932764264418
$ws.Range("L14").Value = 10.16973843113997
$ws.Range("M14").Value = 13.81810107549039
$ws.Range("O14").Value = 25.42905095546805
# Row 15
$ws.Range("B15").Value = 12.161790924279
$ws.Range("C15").Value = 9.06436675562307
$ws.Range("D15").Value = 6.252721549824629
$ws.Range("E15").Value = 12.48597291252311
$ws.Range("F15").Value = 28.48982092971
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 25.54771427247375
$ws.Range("K15").Value = 9.357324523621029
$ws.Range("L15").Value = 10.16964199371849
$ws.Range("M15").Value = 13.81170374187675
$ws.Range("O15").Value = 25.43597626524535
# Row 16
$ws.Range("B16").Value = 11.98174788037809
$ws.Range("C16").Value = 9.039077778526682
$ws.Range("D16").Value = 6.218575159949659
$ws.Range("E16").Value = 12.50201188408262
$ws.Range("F16").Value = 28.50226879049684
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 25.60144486344421
$ws.Range("K16").Value = 9.221126429140341
$ws.Range("L16").Value = 10.16950023317239
$ws.Range("M16").Value = 13.77551695435379
$ws.Range("O16").Value = 25.47715761910229
# Row 17
$ws.Range("B17").Value = 11.87030958624798
$ws.Range("C17").Value = 9.023476028716603
$ws.Range("D17").Value = 6.197407868708395
$ws.Range("E17").Value = 12.51228920799178
$ws.Range("F17").Value = 28.51172511029245
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 25.63559880211361
$ws.Range("K17").Value = 9.136803147433275
$ws.Range("L17").Value = 10.16978067921029
$ws.Range("M17").Value = 13.75374687163814
$ws.Range("O17").Value = 25.50375402768001
# Row 18
$ws.Range("B18").Value = 11.80587282400421
$ws.Range("C18").Value = 9.014470852619938
$ws.Range("D18").Value = 6.185153751258198
$ws.Range("E18").Value = 12.51836148360757
$ws.Range("F18").Value = 28.51783277659011
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 25.65568107526736
$ws.Range("K18").Value = 9.088036206367613
$ws.Range("L18").Value = 10.17007730344543
$ws.Range("M18").Value = 13.74138566125186
$ws.Range("O18").Value = 25.51954095259232
# Row 19
$ws.Range("B19").Value = 11.78399937725227
$ws.Range("C19").Value = 9.011416526619662
$ws.Range("D19").Value = 6.180991244362937
$ws.Range("E19").Value = 12.520445117508
$ws.Range("F19").Value = 28.52001554283624
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 25.6625557659427
$ws.Range("K19").Value = 9.071480417109397
$ws.Range("L19").Value = 10.17020100371162
$ws.Range("M19").Value = 13.73722815248256
$ws.Range("O19").Value = 25.52497013699669
# Row 20
$ws.Range("B20").Value = 11.88220813232184
$ws.Range("C20").Value = 9.025140135814199
$ws.Range("D20").Value = 6.199669398001935
$ws.Range("E20").Value = 12.51117850528051
$ws.Range("F20").Value = 28.51064926164576
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 25.6319177384323
$ws.Range("K20").Value = 9.145807459920762
$ws.Range("L20").Value = 10.16973682676232
$ws.Range("M20").Value = 13.75604779458358
$ws.Range("O20").Value = 25.50087213492155
# Row 21
$ws.Range("B21").Value = 12.20795123940114
$ws.Range("C21").Value = 9.070868083807786
$ws.Range("D21").Value = 6.26146841962907
$ws.Range("E21").Value = 12.48196728169885
$ws.Range("F21").Value = 28.48716988988945
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 25.53421077040009
$ws.Range("K21").Value = 9.392236920741427
$ws.Range("L21").Value = 10.1697919325495
$ws.Range("M21").Value = 13.8211770569957
$ws.Range("O21").Value = 25.42575654703833
# Row 22
$ws.Range("B22").Value = 12.41734428340298
$ws.Range("C22").Value = 9.100456837254292
$ws.Range("D22").Value = 6.3011214399887
$ws.Range("E22").Value = 12.46430805188092
$ws.Range("F22").Value = 28.47778152708348
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 25.47425696134867
$ws.Range("K22").Value = 9.550576240012898
$ws.Range("L22").Value = 10.17102645477342
$ws.Range("M22").Value = 13.8651152415197
$ws.Range("O22").Value = 25.38103224291611
# Row 23
$ws.Range("B23").Value = 12.30592907426315
$ws.Range("C23").Value = 9.084692767484508
$ws.Range("D23").Value = 6.280026696536685
$ws.Range("E23").Value = 12.47360218505496
$ws.Range("F23").Value = 28.48224721892077
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 25.50589806074647
$ws.Range("K23").Value = 9.46633207088693
$ws.Range("L23").Value = 10.17025823363533
$ws.Range("M23").Value = 13.84154177193536
$ws.Range("O23").Value = 25.40450228688652
# Row 24
$ws.Range("B24").Value = 11.87682994491839
$ws.Range("C24").Value = 9.024387903460191
$ws.Range("D24").Value = 6.198647224272744
$ws.Range("E24").Value = 12.51168014403413
$ws.Range("F24").Value = 28.5111335624242
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 25.63358055586033
$ws.Range("K24").Value = 9.141737504510816
$ws.Range("L24").Value = 10.16975623065489
$ws.Range("M24").Value = 13.75500706483619
$ws.Range("O24").Value = 25.50217349281429
# Row 25
$ws.Range("B25").Value = 11.40128735973736
$ws.Range("C25").Value = 8.95814627812927
$ws.Range("D25").Value = 6.107855067381927
$ws.Range("E25").Value = 12.55891101244194
$ws.Range("F25").Value = 28.56784812859129
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 25.78804359965008
$ws.Range("K25").Value = 8.781662443960293
$ws.Range("L25").Value = 10.17439123026892
$ws.Range("M25").Value = 13.69634651108438
$ws.Range("O25").Value = 25.62623061960794

Write-Host "Updated 264 cells in rows 2-25"